$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at K, shifting the table's right border/formatting
# over and giving us a fresh column to hold the new "Tanggal" data.
$ws.Columns("K").Insert()

# The title/period rows (1-2) only ever spanned A:J; the insert incidentally
# stamped a formatted-but-empty K1 cell, so clear it back out.
$ws.Range("K1").Clear()

# Header row (10): insert "Tanggal" ahead of the existing "Status" column by
# relabeling J (now "Tanggal") and giving the new K column the old "Status"
# label - same bold/centered/bordered header style carries over automatically.
$ws.Range("J10").Value = "Tanggal"
$ws.Range("K10").Value = "Status"

# Footer "Tidak ada data" row (11): the boxed placeholder row's right edge
# needs to move from J to the new last column K.
$j11 = $ws.Range("J11")
$j11.Borders.Item(7).LineStyle = -4142
$j11.Borders.Item(10).LineStyle = -4142
$j11.Borders.Item(8).LineStyle = 1
$j11.Borders.Item(9).LineStyle = 1

# Re-merge the placeholder row across the full, now-wider table.
$ws.Range("A11:J11").UnMerge()
$ws.Range("A11:K11").Merge()

# Keep the selection in sync with the newly widened placeholder row.
$ws.Range("A11:K11").Select() | Out-Null
